$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# --- Data corrections on the Property sheet -------------------------------
# Rows 68-74: the "View" column (F) was missing a value; fill it with FALSE
$ws.Range("F68:F74").Value = $false

# Row 75 ("LoadPropertyFinish"): Public/Private/Save flip to FALSE, and the
# View column (F) gets a FALSE value too (previously blank)
$ws.Range("C75:F75").Value = $false

# --- Data validation cleanup -----------------------------------------------
# Replace the old fragmented validation range (F11 F12 F2:F10 F13:F1048576)
# with a single contiguous one, and also apply the same TRUE/FALSE list
# validation to the newly-corrected C75:E75 cells.
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("C75:E75").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# --- Restore the Property sheet as the active / selected tab ---------------
$ws.Activate()
$ws.Range("C75").Select()
